$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to write a value as text, preserving exact string representation
# (prevents Excel from auto-converting numeric-looking strings like "531.20"
# or "1.00" into floating point numbers, which would lose formatting / precision).
function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '71.529.41'
$ws.Range('E2').Value = '  +2.94%  '

# Row 3
Set-TextValue 'D3' '3.997.84'
$ws.Range('E3').Value = '  +1.36%  '

# Row 4
$ws.Range('E4').Value = '  +0.13%  '

# Row 5
Set-TextValue 'D5' '531.20'
$ws.Range('E5').Value = '  +5.43%  '

# Row 6
Set-TextValue 'D6' '149.60'
$ws.Range('E6').Value = '  +1.29%  '

# Row 7
Set-TextValue 'D7' '0.623'

# Row 8
Set-TextValue 'D8' '1.00'
$ws.Range('E8').Value = '  +0.27%  '

# Row 9
Set-TextValue 'D9' '0.738'
$ws.Range('E9').Value = '  +0.44%  '

# Row 10
Set-TextValue 'D10' '0.177'
$ws.Range('E10').Value = '  +0.06%  '

# Row 11
Set-TextValue 'D11' '0.0000345'
$ws.Range('E11').Value = '  -1.50%  '

# Row 12
$ws.Range('E12').Value = '  -0.95%  '

# Row 13
Set-TextValue 'D13' '10.63'
$ws.Range('E13').Value = '  +1.27%  '

# Row 14
Set-TextValue 'D14' '4.645.02'
$ws.Range('E14').Value = '  +1.55%  '

# Row 15
Set-TextValue 'D15' '4.013.31'
$ws.Range('E15').Value = '  +1.58%  '

# Row 16
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D16' '21.25'
$ws.Range('E16').Value = '  +6.15%  '

# Row 17
$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 'D17' '14.31'
$ws.Range('E17').Value = '  +0.51%  '

# Row 18
Set-TextValue 'D18' '1.22'
$ws.Range('E18').Value = '  +2.40%  '

# Row 19
Set-TextValue 'D19' '0.134'
$ws.Range('E19').Value = '  -1.90%  '

# Row 20
Set-TextValue 'D20' '71.497.56'
$ws.Range('E20').Value = '  +2.89%  '

# Row 21
Set-TextValue 'D21' '443.19'
$ws.Range('E21').Value = '  +1.51%  '

# Row 22
Set-TextValue 'D22' '3.55'
$ws.Range('E22').Value = '  +2.58%  '

# Row 23
Set-TextValue 'D23' '91.80'
$ws.Range('E23').Value = '  +3.12%  '

# Row 24
Set-TextValue 'D24' '12.40'
$ws.Range('E24').Value = '  +3.35%  '

# Row 25
Set-TextValue 'D25' '14.27'
$ws.Range('E25').Value = '  -3.11%  '

# Row 26
Set-TextValue 'D26' '4.11'
$ws.Range('E26').Value = '  +5.89%  '

# Row 27
Set-TextValue 'D27' '10.85'
$ws.Range('E27').Value = '  -3.13%  '

# Row 28
Set-TextValue 'D28' '36.99'
$ws.Range('E28').Value = '  -0.53%  '

# Row 29
Set-TextValue 'D29' '13.61'
$ws.Range('E29').Value = '  +1.03%  '

# Row 30
Set-TextValue 'D30' '694.50'
$ws.Range('E30').Value = '  -2.04%  '

# Row 31
$ws.Range('E31').Value = '  -0.15%  '

# Row 32
$ws.Range('E32').Value = '  -0.11%  '

# Row 33
$ws.Range('B33').Value = 'NEARProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D33' '6.78'
$ws.Range('E33').Value = '  +12.08%  '

# Row 34
$ws.Range('B34').Value = 'OKB'
$ws.Range('C34').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D34' '68.33'
$ws.Range('E34').Value = '  +6.04%  '

# Row 35
Set-TextValue 'D35' '0.0₃0911'
$ws.Range('E35').Value = '  +1.76%  '

# Row 36
Set-TextValue 'D36' '0.442'
$ws.Range('E36').Value = '  -2.59%  '

# Row 37
Set-TextValue 'D37' '40.59'
$ws.Range('E37').Value = '  -0.98%  '

# Row 38
Set-TextValue 'D38' '3.55'
$ws.Range('E38').Value = '  +15.10%  '

# Row 39
$ws.Range('E39').Value = '  -0.35%  '

# Row 40
Set-TextValue 'D40' '0.999'
$ws.Range('E40').Value = '  -0.10%  '

# Row 41
Set-TextValue 'D41' '1.00'
$ws.Range('E41').Value = '  -0.11%  '

# Row 42
$ws.Range('E42').Value = '  +0.08%  '

# Row 43
Set-TextValue 'D43' '2.91'
$ws.Range('E43').Value = '  +0.96%  '

# Row 44
$ws.Range('E44').Value = '  +0.98%  '

# Row 45
$ws.Range('E45').Value = '  +11.70%  '

# Row 46
Set-TextValue 'D46' '3.50'
$ws.Range('E46').Value = '  +3.72%  '

# Row 47
$ws.Range('E47').Value = '  +0.21%  '

# Row 48
$ws.Range('B48').Value = 'FLOKI'
$ws.Range('C48').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
Set-TextValue 'D48' '0.000283'
$ws.Range('E48').Value = '  +18.63%  '

# Row 49
$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue 'D49' '9.31'
$ws.Range('E49').Value = '  +5.43%  '

# Row 50
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue 'D50' '0.0₆0358'
$ws.Range('E50').Value = '  +2.74%  '

# Row 51
$ws.Range('E51').Value = '  -0.16%  '
